$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.032.88"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.51%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.358.14"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.97%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.45"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.80%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.70"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.19%  "

# Row 7
$ws.Range("E7").Value = "  +0.14%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.358.24"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.01%  "

# Row 9
$ws.Range("E9").Value = "  -1.65%  "

# Row 10
$ws.Range("E10").Value = "  +0.01%  "

# Row 11
$ws.Range("E11").Value = "  -3.14%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.384"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.77%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.932.37"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.95%  "

# Row 14
$ws.Range("E14").Value = "  +1.32%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.94"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.42%  "

# Row 16
$ws.Range("E16").Value = "  -3.73%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.360.78"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.79%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.164.93"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.48%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.97"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.76%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.80"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.38%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.23"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.86%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "376.46"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.94%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.551"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.45%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.492.61"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.93%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.15%  "

# Row 26
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.06%  "

# Row 27
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000124"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.35%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.77"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +11.29%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.40"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.49%  "

# Row 31
$ws.Range("E31").Value = "  +4.92%  "

# Row 32
$ws.Range("E32").Value = "  -1.98%  "

# Row 33
$ws.Range("E33").Value = "  -0.94%  "

# Row 34
$ws.Range("E34").Value = "  -0.01%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.47"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.13%  "

# Row 36
$ws.Range("E36").Value = "  -5.80%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.75"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.86%  "

# Row 38
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "164.82"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.14%  "

# Row 39
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.52"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.81%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0753"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -5.06%  "

# Row 41
$ws.Range("E41").Value = "  -0.06%  "

# Row 42
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.767"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.70%  "

# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.70"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.97%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.41"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.11%  "

# Row 45
$ws.Range("E45").Value = "  -1.55%  "

# Row 46
$ws.Range("E46").Value = "  -1.65%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.73"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.31%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.15"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.46%  "

# Row 49
$ws.Range("E49").Value = "  -2.65%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.341.74"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.73%  "

# Row 51
$ws.Range("E51").Value = "  -2.08%  "
